# Add a new worksheet "ODI Batting Extra" at the end of the workbook and
# populate it with per-match batting-extras data, mirroring the header
# styling already used by the other data sheets (bold / bordered / centered)
# and keeping MATCH_CODE / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH
# as text values (only BATTING_POSITION is a real number), matching how the
# rest of the workbook stores its "looks numeric" identifiers as text.

$wb = $excel.ActiveWorkbook

# Copy the header formatting from an existing sheet's header row so the new
# header reuses the same cell style rather than creating a brand new one.
$styleSource = $wb.Worksheets.Item("Player Info")
$styleSource.Range("A1:D1").Copy()

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Batting Extra"

$ws.Range("A1:F1").PasteSpecial(-4122)

# Headers
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# MATCH_CODE / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL are stored as text (leading "'"
# keeps Excel from re-interpreting the numeric-looking strings as numbers); blank
# entries are written as a lone "'" so the cell is an empty text cell rather than
# a truly-blank one; BATTING_POSITION is a genuine number.
$data = @(
    @("4284", 2,    "3", "0", "12.06%", "NO"),
    @("4456", $null, $null, $null, $null, "NO"),
    @("4457", 4,    "6", "0", "15.53%", "NO"),
    @("4472", 2,    "8", "0", "47.89%", "NO"),
    @("4473", 2,    "0", "0", $null,    "NO"),
    @("4476", 2,    "0", "0", $null,    "NO"),
    @("4598", $null, $null, $null, $null, "NO"),
    @("4599", $null, $null, $null, $null, "NO"),
    @("4602", 3,    "0", "0", $null,    "NO"),
    @("4660", $null, $null, $null, $null, "NO"),
    @("4663", $null, $null, $null, $null, "NO"),
    @("4666", $null, $null, $null, $null, "NO"),
    @("4698", 2,    "9", "0", "21.77%", "NO"),
    @("4699", 2,    "1", "1", "3.51%",  "NO"),
    @("4700", 2,    "7", "6", "34.10%", "NO"),
    @("4711", 3,    "8", "4", "53.77%", "YES"),
    @("4713", 3,    "0", "0", "3.37%",  "NO"),
    @("4717", 3,    "0", "0", $null,    "NO")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = "'" + $rec[0]

    if ($null -ne $rec[1]) {
        $ws.Cells.Item($row, 2).Value = $rec[1]
    } else {
        $ws.Cells.Item($row, 2).Value = "'"
    }

    if ($null -ne $rec[2]) {
        $ws.Cells.Item($row, 3).Value = "'" + $rec[2]
    } else {
        $ws.Cells.Item($row, 3).Value = "'"
    }

    if ($null -ne $rec[3]) {
        $ws.Cells.Item($row, 4).Value = "'" + $rec[3]
    } else {
        $ws.Cells.Item($row, 4).Value = "'"
    }

    if ($null -ne $rec[4]) {
        $ws.Cells.Item($row, 5).Value = "'" + $rec[4]
    } else {
        $ws.Cells.Item($row, 5).Value = "'"
    }

    $ws.Cells.Item($row, 6).Value = $rec[5]

    $row++
}

# Restore the original active sheet/selection so this edit only adds the new
# sheet without disturbing the workbook's view state.
$wb.Worksheets.Item(1).Activate()
$null = $wb.Worksheets.Item(1).Range("A1").Select()

